$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the two data rows (c09d2861... and 4d375c79...) swap
# identity - row 2 becomes the c09d2861 file, row 3 becomes the 4d375c79
# file - and the 4d375c79 row picks up a new "Ready for handoff" status and
# timestamp (a fresh handoff report was generated).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.md'
$wsOverview.Range("A3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.md'

$wsOverview.Range("E3").Value = 'Ready for handoff'
$wsOverview.Range("F3").Value = 'Ready for handoff'
$wsOverview.Range("G3").Value = '2016-09-05 10:54:57'

# Hyperlinks on B2/B3 need their display text swapped while keeping the same
# underlying targets, so drop and recreate them in the desired order.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", 'e2e\c09d2861-933d-4f7e-a942-dfe71323448b.md')
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", 'e2e\4d375c79-a9db-48a3-b453-a6097fcdfc18.md')

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus the 4d375c79 row's status/date/error.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.md'
$wsZhCn.Range("G2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.zh-cn.xlf'
$wsZhCn.Range("I2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.md'
$wsZhCn.Range("J2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.zh-cn.xlf'

$wsZhCn.Range("A3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.md'
$wsZhCn.Range("C3").Value = 'Ready for handoff'
$wsZhCn.Range("G3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.zh-cn.xlf'
$wsZhCn.Range("H3").Value = '2016-09-05 10:54:45'
$wsZhCn.Range("I3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.md'
$wsZhCn.Range("J3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.zh-cn.xlf'
$wsZhCn.Range("P3").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53899237570685e41dfe068e60bdf6b0af27382f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md.'

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", 'c09d2861-933d-4f7e-a942-dfe71323448b.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d651fd4ae8d0f449432bd188f67ea762fb1320d5/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", 'c09d2861-933d-4f7e-a942-dfe71323448b.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", '4d375c79-a9db-48a3-b453-a6097fcdfc18.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d651fd4ae8d0f449432bd188f67ea762fb1320d5/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", '4d375c79-a9db-48a3-b453-a6097fcdfc18.md')

# Error Detail column got wide enough to actually show the message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Sheet "de-de": same pattern as zh-cn, different xlf/url suffixes.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.md'
$wsDeDe.Range("G2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.de-de.xlf'
$wsDeDe.Range("I2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.md'
$wsDeDe.Range("J2").Value = 'c09d2861-933d-4f7e-a942-dfe71323448b.641a8eac3ce33e7ae713ddd97cca8833e30ce46f.de-de.xlf'

$wsDeDe.Range("A3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.md'
$wsDeDe.Range("C3").Value = 'Ready for handoff'
$wsDeDe.Range("G3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.de-de.xlf'
$wsDeDe.Range("H3").Value = '2016-09-05 10:54:57'
$wsDeDe.Range("I3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.md'
$wsDeDe.Range("J3").Value = '4d375c79-a9db-48a3-b453-a6097fcdfc18.893ac78b6598f9ffc79a2cd242b48e478a4d91e3.de-de.xlf'
$wsDeDe.Range("P3").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53899237570685e41dfe068e60bdf6b0af27382f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md.'

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", 'c09d2861-933d-4f7e-a942-dfe71323448b.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3a1a13cabe2b52d682ceb57c7caca4a61b2565aa/e2e/4d375c79-a9db-48a3-b453-a6097fcdfc18.md", "", "", 'c09d2861-933d-4f7e-a942-dfe71323448b.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b02cbc8fa2c30700d268c53b7006b3550719c9f/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", '4d375c79-a9db-48a3-b453-a6097fcdfc18.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3a1a13cabe2b52d682ceb57c7caca4a61b2565aa/e2e/c09d2861-933d-4f7e-a942-dfe71323448b.md", "", "", '4d375c79-a9db-48a3-b453-a6097fcdfc18.md')

$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
